$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Mann-Whitney results: p-values (alpha_MW) and significance flags
# recomputed for a newer run of the test; the Nutrient column (A) is unchanged.

$ws.Range("B2").Value = 0.000055527508301138
$ws.Range("C2").Value = "yes"

$ws.Range("B3").Value = 0.00000000000681583725823273
$ws.Range("C3").Value = "yes"

$ws.Range("B4").Value = 0.660538202472032
$ws.Range("C4").Value = "no"

$ws.Range("B5").Value = 0.0114055311330627
$ws.Range("C5").Value = "yes"

$ws.Range("B6").Value = 0.462343259588609
$ws.Range("C6").Value = "no"

$ws.Range("B7").Value = 0.0000000000481337518394799
$ws.Range("C7").Value = "yes"

$ws.Range("B8").Value = 0.0000000000000000750761848812129
$ws.Range("C8").Value = "yes"

$ws.Range("B9").Value = 0.278465241197607
$ws.Range("C9").Value = "no"
